$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.281.06'
$ws.Range('E2').Value = '  -2.28%  '
$ws.Range('D3').Value = '1.878.54'
$ws.Range('E3').Value = '  -1.74%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '236.25'
$ws.Range('E5').Value = '  -1.25%  '
$ws.Range('E6').Value = '  +0.16%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4854'
$ws.Range('E7').Value = '  -1.11%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2872'
$ws.Range('E8').Value = '  -3.21%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06585'
$ws.Range('E9').Value = '  -3.01%  '
$ws.Range('D10').Value = '1.881.44'
$ws.Range('E10').Value = '  -1.65%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '16.73'
$ws.Range('E11').Value = '  -2.25%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.190'
$ws.Range('E13').Value = '  +1.21%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '87.11'
$ws.Range('E14').Value = '  -3.35%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6556'
$ws.Range('E15').Value = '  -2.44%  '
$ws.Range('D16').Value = '30.243.59'
$ws.Range('E16').Value = '  -2.22%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.34'
$ws.Range('E17').Value = '  -1.30%  '
$ws.Range('E18').Value = '  +0.20%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007726'
$ws.Range('E19').Value = '  -3.08%  '
$ws.Range('D20').Value = '2.149.50'
$ws.Range('E20').Value = '  +0.28%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.294'
$ws.Range('E21').Value = '  +3.67%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.000'
$ws.Range('E22').Value = '  +0.16%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '194.86'
$ws.Range('E23').Value = '  -6.51%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.129'
$ws.Range('E24').Value = '  -1.50%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.281'
$ws.Range('E25').Value = '  -4.00%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '159.05'
$ws.Range('E26').Value = '  +0.55%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.07'
$ws.Range('E27').Value = '  -4.88%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.915'
$ws.Range('E28').Value = '  -2.66%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.440'
$ws.Range('E29').Value = '  +0.98%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.267'
$ws.Range('E30').Value = '  -1.46%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09120'
$ws.Range('E31').Value = '  -0.66%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.054'
$ws.Range('E32').Value = '  +0.37%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05091'
$ws.Range('E33').Value = '  -1.80%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7178'
$ws.Range('E34').Value = '  -4.52%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.094'
$ws.Range('E35').Value = '  -2.32%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.712'
$ws.Range('E36').Value = '  +0.82%  '
$ws.Range('E37').Value = '  -2.66%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.636'
$ws.Range('E38').Value = '  -3.55%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.9197'
$ws.Range('E39').Value = '  -0.95%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.039'
$ws.Range('E40').Value = '  -3.78%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '106.30'
$ws.Range('E41').Value = '  -0.31%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4274'
$ws.Range('E42').Value = '  -4.97%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.787'
$ws.Range('E43').Value = '  -0.74%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.000'
$ws.Range('E44').Value = '  -0.63%  '
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.384'
$ws.Range('E45').Value = '  -5.19%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '65.92'
$ws.Range('E46').Value = '  -1.29%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1319'
$ws.Range('E47').Value = '  -4.26%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.928'
$ws.Range('E48').Value = '  -0.41%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05748'
$ws.Range('E49').Value = '  -2.63%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '33.84'
$ws.Range('E50').Value = '  -3.69%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3815'
$ws.Range('E51').Value = '  -6.46%  '
